# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps on the per-language
# status sheets, as produced by a re-run of the handback report
# generation.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 06:00:51"
$wsZhCn.Range("H2").Value = "2016-03-18 06:01:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 06:00:54"
$wsDeDe.Range("H2").Value = "2016-03-18 06:01:15"
